# Estado de cuenta: se eliminan los periodos de mora anteriores y se
# reescriben con los nuevos datos (orden descendente por Periodo Mora),
# dejando la fila de WILMER (NIT 79727775 / periodo 1902) al final de la
# tabla. Los formatos de cada fila (bordes, relleno, etc.) permanecen en
# su posicion original; solo se actualizan los valores.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$docNum   = "1051357849"
$docName  = "ZILLAH DE JESUS PATERNOSTRO CASTILLO"
$periods  = @("2504","2503","2502","2501","2412","2411","2410","2409","2408","2407","2406","2405","2404","2403","2402","2401","2312","2311","2310","2309","2308","2307","2306","2305","2304","2303")

$row = 16
foreach ($p in $periods) {
    $ws.Cells.Item($row, 3).Value = $docNum
    $ws.Cells.Item($row, 4).Value = $docName
    $ws.Cells.Item($row, 5).Value = $p
    if ($row -eq 16) {
        $ws.Cells.Item($row, 6).Value = 62866
    } else {
        $ws.Cells.Item($row, 6).Value = 82000
    }
    $ws.Cells.Item($row, 7).Value = 2050000
    $row = $row + 1
}

# Ultima fila de la tabla: trabajador WILMER al final, con sus valores originales.
$ws.Cells.Item(42, 3).Value = "79727775"
$ws.Cells.Item(42, 4).Value = "WILMER DINAEL GARCIA MAYORCA"
$ws.Cells.Item(42, 5).Value = "1902"
$ws.Cells.Item(42, 6).Value = 112000
$ws.Cells.Item(42, 7).Value = 3000000
